$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Sheet1 -> Sheet2
$ws.Name = "Sheet2"

# Update the active selection to D15
$ws.Range("D15").Select()
